# Update the "F1_Venta_23_Ene_Porcentaje" percentage figures (column C,
# rows 2-13) and move the active selection, matching the uploaded sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C2").Value  = 1.01
$ws.Range("C3").Value  = 0.51
$ws.Range("C4").Value  = 0.67
$ws.Range("C5").Value  = 0.47
$ws.Range("C6").Value  = 0.87
$ws.Range("C7").Value  = 0.82
$ws.Range("C8").Value  = 0.73
$ws.Range("C9").Value  = 0.68
$ws.Range("C10").Value = 0.86
$ws.Range("C11").Value = 0.71
$ws.Range("C12").Value = 0.66
$ws.Range("C13").Value = 0.42

# Reflect the saved cursor position (the last thing the author clicked
# before uploading).
$ws.Range("A14").Select()
